# Updates to column C (count) and column D (amount) for the
# "Fonds de solidarite" regional dataset, refreshed with 2020-05-24 data.
# Each entry is: row number, new C value, new D value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @(2, 23877, 34851567),
    @(3, 60223, 89123815),
    @(4, 20440, 30493309),
    @(5, 5388, 8058820),
    @(6, 1073, 1607197),
    @(10, 25604, 35364098),
    @(11, 6269, 9169063),
    @(12, 17677, 26137883),
    @(13, 5509, 8225163),
    @(14, 1305, 1952095),
    @(15, 236, 351266),
    @(17, 6424, 8691614),
    @(18, 8710, 12695451),
    @(19, 21621, 31997208),
    @(20, 6874, 10269407),
    @(21, 1650, 2469802),
    @(24, 7472, 10211253),
    @(25, 4940, 7205853),
    @(26, 15237, 22538519),
    @(27, 5192, 7759432),
    @(28, 1250, 1874491),
    @(31, 5287, 7096203),
    @(32, 1744, 2513803),
    @(33, 4600, 6755754),
    @(34, 1873, 2786851),
    @(35, 487, 727041),
    @(38, 1181, 1607948),
    @(39, 11137, 16233908),
    @(40, 34330, 50737254),
    @(41, 12742, 19033677),
    @(42, 3533, 5288416),
    @(43, 615, 921436),
    @(46, 10435, 14205878),
    @(47, 1014, 1467198),
    @(48, 3787, 5584834),
    @(49, 1409, 2106464),
    @(50, 438, 654500),
    @(52, 2503, 3482429),
    @(54, 963, 1426244),
    @(59, 469, 666625),
    @(60, 10069, 14625136),
    @(61, 30695, 45282625),
    @(62, 10677, 15954655),
    @(63, 2960, 4430138),
    @(64, 522, 781931),
    @(67, 9978, 13368699),
    @(68, 2746, 4005992),
    @(69, 7449, 10989910),
    @(70, 2641, 3945983),
    @(71, 868, 1300049),
    @(74, 2903, 3956945),
    @(75, 881, 1278393),
    @(76, 3015, 4453499),
    @(77, 1201, 1796939),
    @(79, 85, 127069),
    @(81, 1797, 2413902),
    @(83, 109, 163110),
    @(87, 7079, 10353049),
    @(88, 20363, 30122799),
    @(89, 6696, 10007887),
    @(90, 1769, 2648271),
    @(94, 6356, 8570116),
    @(95, 19448, 28242536),
    @(96, 45186, 66645947),
    @(97, 14471, 21604981),
    @(98, 3873, 5794340),
    @(99, 661, 989862),
    @(102, 16646, 22609010),
    @(103, 22224, 32313700),
    @(104, 50288, 74090263),
    @(105, 15749, 23482059),
    @(106, 4021, 6007501),
    @(107, 656, 981054),
    @(109, 5, 7500),
    @(110, 19712, 26585844),
    @(111, 8643, 12623139),
    @(112, 22405, 33157573),
    @(113, 7794, 11635128),
    @(114, 1885, 2819669),
    @(118, 7068, 9618736),
    @(119, 21653, 31479362),
    @(120, 53421, 78781854),
    @(121, 16119, 24072780),
    @(122, 4005, 5990997),
    @(123, 812, 1216212),
    @(126, 18424, 25314565),
    @(127, 29809, 43654042),
    @(128, 89658, 132797538),
    @(129, 39824, 59518407),
    @(130, 12413, 18591962),
    @(131, 2533, 3793255),
    @(135, 29310, 40855546)
)

foreach ($u in $updates) {
    $row = $u[0]
    $cVal = $u[1]
    $dVal = $u[2]
    $ws.Cells.Item($row, 3).Value = $cVal
    $ws.Cells.Item($row, 4).Value = $dVal
}
